# Check-list.xlsx — 09.05.2023 Глумов v1.2 final
# ---------------------------------------------------------------
# 1) Sheet "Чеклист": mark item 6 (row 7) as done -> B7 = "+"
# 2) Sheet "План": mark several plan items as done in column C,
#    and drop the now-obsolete note about the free course (B37).
# ---------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- Sheet: Чеклист -------------------------------------------------
$checklist = $wb.Worksheets("Чеклист")
$checklist.Range("B7").Value = "+"
$checklist.Range("B8").Select()

# --- Sheet: План ------------------------------------------------------
$plan = $wb.Worksheets("План")

# Этап 2 — оба пункта выполнены
$plan.Range("C33").Value = "Выполнено"
$plan.Range("C34").Value = "Выполнено"

# Заметка про наработки из бесплатного курса больше не актуальна
$plan.Range("B37").Clear()
$plan.Rows(37).RowHeight = 45

# Этап 3 — оба пункта выполнены (строчными буквами)
$plan.Range("C38").Value = "выполнено"
$plan.Range("C40").Value = "выполнено"

$plan.Activate()
$plan.Range("C29").Select()
